$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 275.46667
$ws.Range("I33").Value = 176.36363
$ws.Range("K33").Value = 176.36363
$ws.Range("M33").Value = 52.63637
$ws.Range("H40").Value = 5923
$ws.Range("I40").Value = 5186.125
$ws.Range("K40").Value = 5186.125
$ws.Range("M40").Value = -5011.125
$ws.Range("H53").Value = 858.38464
$ws.Range("I53").Value = 1094.8
$ws.Range("K53").Value = 1094.8
$ws.Range("M53").Value = -457.8
$ws.Range("H58").Value = 1501.3572
$ws.Range("J58").Value = 6500
$ws.Range("L58").Value = 19500
$ws.Range("N58").Value = -19800
$ws.Range("H125").Value = 949.35297
$ws.Range("I125").Value = 932.3333
$ws.Range("J125").Value = 958.63635
$ws.Range("K125").Value = 8390.9997
$ws.Range("L125").Value = 8627.727150000001
$ws.Range("M125").Value = -5930.9997
$ws.Range("N125").Value = -13547.72715
$ws.Range("H133").Value = 93994.336
$ws.Range("J133").Value = 93994.336
$ws.Range("L133").Value = 93994.336
$ws.Range("N133").Value = -104114.336
$ws.Range("H134").Value = 71102.8
$ws.Range("J134").Value = 71102.8
$ws.Range("L134").Value = 71102.8
$ws.Range("N134").Value = -81242.8
$ws.Range("H135").Value = 610.61536
$ws.Range("I135").Value = 610.61536
$ws.Range("K135").Value = 5495.53824
$ws.Range("M135").Value = -2960.53824
$ws.Range("H136").Value = 105994.5
$ws.Range("J136").Value = 105994.5
$ws.Range("L136").Value = 105994.5
$ws.Range("N136").Value = -116194.5
$ws.Range("H138").Value = 305526.75
$ws.Range("I138").Value = 3329.7693
$ws.Range("J138").Value = 361649.03
$ws.Range("K138").Value = 9989.3079
$ws.Range("L138").Value = 1084947.09
$ws.Range("M138").Value = -4849.3079
$ws.Range("N138").Value = -1095227.09
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13327.431
$ws.Range("I32").Value = 8432.078
$ws.Range("J32").Value = 34214.266
$ws.Range("K32").Value = 8432.078
$ws.Range("L32").Value = 34214.266
$ws.Range("M32").Value = -8145.078
$ws.Range("N32").Value = -34788.266
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716
$ws.Range("H122").Value = 5474.8945
$ws.Range("I122").Value = 3091.125
$ws.Range("K122").Value = 9273.375
$ws.Range("M122").Value = -6823.375
$ws.Range("H132").Value = 2678.0334
$ws.Range("I132").Value = 1944.7778
$ws.Range("J132").Value = 3777.9167
$ws.Range("K132").Value = 5834.3334
$ws.Range("L132").Value = 11333.7501
$ws.Range("M132").Value = -3304.3334
$ws.Range("N132").Value = -16393.7501
$ws.Range("H133").Value = 113577
$ws.Range("J133").Value = 113577
$ws.Range("L133").Value = 113577
$ws.Range("N133").Value = -118637
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 2999
$ws.Range("J19").Value = 2999
$ws.Range("L19").Value = 2999
$ws.Range("N19").Value = -3345
$ws.Range("H35").Value = 38713.668
$ws.Range("J35").Value = 59074
$ws.Range("L35").Value = 59074
$ws.Range("N35").Value = -59694
$ws.Range("H59").Value = 128998.664
$ws.Range("J59").Value = 128998.664
$ws.Range("L59").Value = 128998.664
$ws.Range("N59").Value = -130692.664
$ws.Range("H82").Value = 85455.8
$ws.Range("I82").Value = 72498.5
$ws.Range("K82").Value = 72498.5
$ws.Range("M82").Value = -72115.5
$ws.Range("H85").Value = 85455.8
$ws.Range("I85").Value = 72498.5
$ws.Range("K85").Value = 72498.5
$ws.Range("M85").Value = -71172.5
$ws.Range("H94").Value = 100000840
$ws.Range("I94").Value = 111111940
$ws.Range("J94").Value = 998
$ws.Range("K94").Value = 111111940
$ws.Range("L94").Value = 998
$ws.Range("M94").Value = -111111489
$ws.Range("N94").Value = -1900
$ws.Range("H105").Value = 26003300
$ws.Range("I105").Value = 2502223.8
$ws.Range("J105").Value = 41670684
$ws.Range("K105").Value = 2502223.8
$ws.Range("L105").Value = 41670684
$ws.Range("M105").Value = -2500476.8
$ws.Range("N105").Value = -41674178
$ws.Range("H108").Value = 85000
$ws.Range("I108").Value = 85000
$ws.Range("K108").Value = 85000
$ws.Range("M108").Value = -81160
$ws.Range("H134").Value = 3890
$ws.Range("I134").Value = 3392.2144
$ws.Range("J134").Value = 4885.5713
$ws.Range("K134").Value = 10176.6432
$ws.Range("L134").Value = 14656.7139
$ws.Range("M134").Value = -7641.643199999999
$ws.Range("N134").Value = -19726.7139
$ws.Range("H139").Value = 104989
$ws.Range("J139").Value = 104989
$ws.Range("L139").Value = 104989
$ws.Range("N139").Value = -115269
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1990
$ws.Range("I22").Value = 1225
$ws.Range("K22").Value = 1225
$ws.Range("M22").Value = -875
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 3092564.2
$ws.Range("I31").Value = 5100.3076
$ws.Range("J31").Value = 15635386
$ws.Range("K31").Value = 5100.3076
$ws.Range("L31").Value = 15635386
$ws.Range("M31").Value = -4805.3076
$ws.Range("N31").Value = -15635976
$ws.Range("H34").Value = 3092564.2
$ws.Range("I34").Value = 5100.3076
$ws.Range("J34").Value = 15635386
$ws.Range("K34").Value = 5100.3076
$ws.Range("L34").Value = 15635386
$ws.Range("M34").Value = -4898.3076
$ws.Range("N34").Value = -15635790
$ws.Range("H105").Value = 2111.4666
$ws.Range("I105").Value = 1181.1111
$ws.Range("K105").Value = 1181.1111
$ws.Range("M105").Value = 565.8888999999999
$ws.Range("H132").Value = 10003389
$ws.Range("I132").Value = 10872205
$ws.Range("K132").Value = 32616615
$ws.Range("M132").Value = -32614085
$ws.Range("H134").Value = 3929.7
$ws.Range("I134").Value = 2197.7036
$ws.Range("K134").Value = 6593.110799999999
$ws.Range("M134").Value = -4058.110799999999
$ws.Range("H141").Value = 574029.25
$ws.Range("J141").Value = 574029.25
$ws.Range("L141").Value = 574029.25
$ws.Range("N141").Value = -584389.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 994.5
$ws.Range("I2").Value = 826.1429000000001
$ws.Range("J2").Value = 1101.6364
$ws.Range("K2").Value = 4956.857400000001
$ws.Range("L2").Value = 6609.8184
$ws.Range("M2").Value = -4843.857400000001
$ws.Range("N2").Value = -6835.8184
$ws.Range("H13").Value = 133
$ws.Range("I13").Value = 149.5
$ws.Range("K13").Value = 448.5
$ws.Range("M13").Value = -280.5
$ws.Range("H16").Value = 1800
$ws.Range("J16").Value = 1800
$ws.Range("L16").Value = 5400
$ws.Range("N16").Value = -5746
$ws.Range("H38").Value = 539.36365
$ws.Range("I38").Value = 320.4
$ws.Range("J38").Value = 721.8333
$ws.Range("K38").Value = 961.1999999999999
$ws.Range("L38").Value = 2165.4999
$ws.Range("M38").Value = -614.1999999999999
$ws.Range("N38").Value = -2859.4999
$ws.Range("H80").Value = 8666
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 8666
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H122").Value = 1450.037
$ws.Range("I122").Value = 1345.125
$ws.Range("J122").Value = 1494.2106
$ws.Range("K122").Value = 12106.125
$ws.Range("L122").Value = 13447.8954
$ws.Range("M122").Value = -9656.125
$ws.Range("N122").Value = -18347.8954
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2750.5264
$ws.Range("I132").Value = 2623
$ws.Range("J132").Value = 2824.9167
$ws.Range("K132").Value = 7869
$ws.Range("L132").Value = 8474.750100000001
$ws.Range("M132").Value = -5339
$ws.Range("N132").Value = -13534.7501
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 78642.336
$ws.Range("J140").Value = 78642.336
$ws.Range("L140").Value = 78642.336
$ws.Range("N140").Value = -89002.336
$ws.Range("H141").Value = 113549.555
$ws.Range("J141").Value = 113549.555
$ws.Range("L141").Value = 113549.555
$ws.Range("N141").Value = -123909.555
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 24999.777
$ws.Range("J52").Value = 24999.875
$ws.Range("L52").Value = 24999.875
$ws.Range("N52").Value = -25451.875
$ws.Range("H135").Value = 57713.43
$ws.Range("J135").Value = 57713.43
$ws.Range("L135").Value = 57713.43
$ws.Range("N135").Value = -67853.42999999999
$ws.Range("H137").Value = 76998.13
$ws.Range("J137").Value = 76998.13
$ws.Range("L137").Value = 76998.13
$ws.Range("N137").Value = -87198.13
$ws.Range("H139").Value = 79999
$ws.Range("J139").Value = 79999
$ws.Range("L139").Value = 79999
$ws.Range("N139").Value = -90279
$ws.Range("H141").Value = 69999
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359
